# Update cryptocurrency price (D) and 1h volume-change (E) columns
# with refreshed data, matching the upstream GitHub Actions scrape run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.805.80"
$ws.Range("E2").Value = "  -0.22%  "
$ws.Range("D3").Value = "1.635.64"
$ws.Range("E3").Value = "  -0.21%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'215.24"
$ws.Range("E5").Value = "  -0.34%  "
$ws.Range("E6").Value = "  -0.64%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").Value = "'0.0642"
$ws.Range("E9").Value = "  -0.29%  "
$ws.Range("D10").Value = "'19.86"
$ws.Range("E10").Value = "  +1.10%  "
$ws.Range("D11").Value = "'0.0785"
$ws.Range("E11").Value = "  +0.96%  "
$ws.Range("E12").Value = "  -0.92%  "
$ws.Range("D13").Value = "1.643.18"
$ws.Range("E13").Value = "  +0.48%  "
$ws.Range("D14").Value = "1.861.16"
$ws.Range("E14").Value = "  -0.16%  "
$ws.Range("D15").Value = "'0.557"
$ws.Range("E15").Value = "  -1.28%  "
$ws.Range("E16").Value = "  +1.55%  "
$ws.Range("D17").Value = "'63.07"
$ws.Range("E17").Value = "  -0.34%  "
$ws.Range("D18").Value = "25.818.56"
$ws.Range("E18").Value = "  -0.23%  "
$ws.Range("E19").Value = "  -0.17%  "
$ws.Range("E20").Value = "  +2.39%  "
$ws.Range("D21").Value = "'194.04"
$ws.Range("E21").Value = "  -0.40%  "
$ws.Range("D22").Value = "'9.94"
$ws.Range("E22").Value = "  +0.34%  "
$ws.Range("E23").Value = "  +0.48%  "
$ws.Range("E24").Value = "  -0.07%  "
$ws.Range("D25").Value = "'1.77"
$ws.Range("E25").Value = "  -0.17%  "
$ws.Range("D26").Value = "'139.15"
$ws.Range("E26").Value = "  -0.77%  "
$ws.Range("D27").Value = "'0.121"
$ws.Range("E27").Value = "  -5.03%  "
$ws.Range("D28").Value = "'6.84"
$ws.Range("E28").Value = "  +0.85%  "
$ws.Range("D29").Value = "'15.58"
$ws.Range("E29").Value = "  +0.79%  "
$ws.Range("E30").Value = "  -0.01%  "
$ws.Range("D31").Value = "'0.0498"
$ws.Range("E31").Value = "  +1.91%  "
$ws.Range("D32").Value = "'3.34"
$ws.Range("E32").Value = "  +1.29%  "
$ws.Range("E33").Value = "  +1.29%  "
$ws.Range("E34").Value = "  +2.15%  "
$ws.Range("E35").Value = "  +0.67%  "
$ws.Range("D36").Value = "'0.899"
$ws.Range("E36").Value = "  -0.84%  "
$ws.Range("D37").Value = "'2.58"
$ws.Range("E37").Value = "  -0.13%  "
$ws.Range("D38").Value = "'0.551"
$ws.Range("E38").Value = "  +0.11%  "
$ws.Range("D39").Value = "1.106.40"
$ws.Range("E39").Value = "  -2.07%  "
$ws.Range("D40").Value = "'0.0156"
$ws.Range("E40").Value = "  -0.06%  "
$ws.Range("E41").Value = "  +0.58%  "
$ws.Range("E42").Value = "  +0.77%  "
$ws.Range("D43").Value = "'0.803"
$ws.Range("E43").Value = "  +0.25%  "
$ws.Range("D44").Value = "'99.17"
$ws.Range("E44").Value = "  +1.52%  "
$ws.Range("D45").Value = "0.0₆0107"
$ws.Range("E45").Value = "  -5.54%  "
$ws.Range("D46").Value = "'55.60"
$ws.Range("E46").Value = "  +0.09%  "
$ws.Range("D47").Value = "'2.50"
$ws.Range("E47").Value = "  +12.29%  "
$ws.Range("E48").Value = "  -6.03%  "
$ws.Range("E49").Value = "  -0.36%  "
$ws.Range("E50").Value = "  -0.41%  "
$ws.Range("E51").Value = "  -0.15%  "
